# mppsteel Feedstock Prices.xlsx - remove "Excel Tab" column (G) data,
# select cell H5, and set zoom to 85%.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column G (the "Excel Tab" header + "Price & Emission parameters" values)
$ws.Range("G1:G7").ClearContents()

# Update the view: zoom to 85% and select H5 (matches the author's saved view state)
[void]$ws.Range("H5").Select()
$excel.ActiveWindow.Zoom = 85
